# Update countries & provincias Spain
# The underlying data table (Pais sheet) is kept sorted by column B
# ("Casos totales") descending. The source data refreshed a handful of
# per-country totals; those updates shifted several countries' rank,
# so both the country label (column A) and the stats (columns B:H) for
# the affected rows need to be rewritten to match the new sort order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Country, $Total, $New, $Active, $Recovered, $Critical, $DeathsToday, $Deaths) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $New
    $ws.Cells.Item($Row, 4).Value = $Active
    $ws.Cells.Item($Row, 5).Value = $Recovered
    $ws.Cells.Item($Row, 6).Value = $Critical
    $ws.Cells.Item($Row, 7).Value = $DeathsToday
    $ws.Cells.Item($Row, 8).Value = $Deaths
}

Set-Row 14  "Iran"                  143849 2258 112988 23234 0 63 7627
Set-Row 26  "Bielorrusia"           39858  902  16660  22979 0 5  219
Set-Row 30  "Emiratos Arabes Unidos" 32532 563  16685  15589 0 3  258
Set-Row 32  "Suiza"                 30796  20   28300  579   0 0  1917
Set-Row 36  "Kuwait"                24112  845  8698   15229 0 10 185
Set-Row 37  "Colombia"              24104  0    6111   17190 0 0  803
Set-Row 41  "Rumania"               18791  197  12629  4933  0 2  1229
Set-Row 56  "Oman"                  9009   636  2177   6792  0 1  40
Set-Row 57  "Argelia"               8857   0    5129   3105  0 0  623
Set-Row 58  "Nigeria"               8733   0    2501   5978  0 0  254
Set-Row 59  "Noruega"               8401   0    7727   439   0 0  235
Set-Row 78  "Senegal"               3348   95   1686   1623  0 1  39
Set-Row 79  "Guinea"                3275   0    1673   1582  0 0  20
Set-Row 132 "San Marino"            670    3    322    306   0 0  42
Set-Row 197 "Fiyi"                  18     0    15     3     0 0  0
Set-Row 198 "Curazao"               18     0    14     3     0 0  1
Set-Row 199 "Nueva Caledonia"       18     0    18     0     0 0  0
Set-Row 200 "Belice"                18     0    16     0     0 0  2
Set-Row 201 "Santa Lucia"           18     0    18     0     0 0  0

# Timestamp banner in row 1
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 12:25"
